# New PO forecast model
# Updates three sheets: "Weekly Quantity", "Monthly Trend", "PO Forecast"

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# -----------------------------------------------------------------
# Sheet "Weekly Quantity": append row 32
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Cells.Item(32, 1).Value = 45662.99999999999
$ws1.Cells.Item(32, 1).NumberFormat = $dateFmt
$ws1.Cells.Item(32, 2).Value = 50

# -----------------------------------------------------------------
# Sheet "Monthly Trend": append row 15
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(15, 1).Value = 45688.99999999999
$ws2.Cells.Item(15, 1).NumberFormat = $dateFmt
$ws2.Cells.Item(15, 2).Value = 50

# -----------------------------------------------------------------
# Sheet "PO Forecast": refreshed forecast numbers + one more week
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

# Shift the last stretch of forecast dates forward (rows 32-39) and
# add a brand-new row 40.
$ws3.Cells.Item(32, 1).Value = 45662.99999999999
$ws3.Cells.Item(33, 1).Value = 45669.99999999999
$ws3.Cells.Item(34, 1).Value = 45676.99999999999
$ws3.Cells.Item(35, 1).Value = 45683.99999999999
$ws3.Cells.Item(36, 1).Value = 45690.99999999999
$ws3.Cells.Item(37, 1).Value = 45697.99999999999
$ws3.Cells.Item(38, 1).Value = 45704.99999999999
$ws3.Cells.Item(39, 1).Value = 45711.99999999999

$ws3.Cells.Item(40, 1).Value = 45718.99999999999
$ws3.Cells.Item(40, 1).NumberFormat = $dateFmt
$ws3.Cells.Item(40, 2).Value = 65

# Revised forecast quantities (PO_Forecast column)
$ws3.Cells.Item(2, 2).Value = 35
$ws3.Cells.Item(4, 2).Value = 37
$ws3.Cells.Item(5, 2).Value = 38
$ws3.Cells.Item(14, 2).Value = 53
$ws3.Cells.Item(16, 2).Value = 54
$ws3.Cells.Item(17, 2).Value = 54
$ws3.Cells.Item(19, 2).Value = 55
$ws3.Cells.Item(20, 2).Value = 55
$ws3.Cells.Item(21, 2).Value = 56
$ws3.Cells.Item(22, 2).Value = 57
$ws3.Cells.Item(23, 2).Value = 57
$ws3.Cells.Item(24, 2).Value = 58
$ws3.Cells.Item(25, 2).Value = 58
$ws3.Cells.Item(26, 2).Value = 58
$ws3.Cells.Item(27, 2).Value = 59
$ws3.Cells.Item(28, 2).Value = 59
$ws3.Cells.Item(29, 2).Value = 60
$ws3.Cells.Item(30, 2).Value = 60
$ws3.Cells.Item(31, 2).Value = 61
